# Fruta / hortaliza, semanal
# Insert a new weekly block of 4 rows (one per "Variedad") at the top of the
# data table (just above the existing first block, which started at row 959),
# pushing all subsequent rows down by 4. The new block uses Fecha = 44585
# (2022-01-24) and its own Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 959; everything from 959 downward shifts
# down to 963 onward (so the last existing row, 1018, ends up at 1022).
$ws.Rows("959:962").Insert()

# Columns that are constant across every data row in this sheet.
$commonA = 5
$commonB = "Macroferia Regional de Talca"
$commonC = "Maule"
$commonE = 7
$commonF = 100112033
$commonG = "Lechuga"
$commonR = "Hortaliza"

# Row 959 - Conconina(o), Primera
$ws.Range("A959").Value = $commonA
$ws.Range("B959").Value = $commonB
$ws.Range("C959").Value = $commonC
$ws.Range("D959").Value = 44585
$ws.Range("E959").Value = $commonE
$ws.Range("F959").Value = $commonF
$ws.Range("G959").Value = $commonG
$ws.Range("H959").Value = "Conconina(o)"
$ws.Range("I959").Value = "Primera"
$ws.Range("J959").Value = 400
$ws.Range("K959").Value = 4500
$ws.Range("L959").Value = 4500
$ws.Range("M959").Value = 4500
$ws.Range("N959").Value = "$/caja 10 unidades"
$ws.Range("O959").Value = "Región del Maule"
$ws.Range("P959").Value = 450
$ws.Range("Q959").Value = 10
$ws.Range("R959").Value = $commonR

# Row 960 - Escarola, Primera
$ws.Range("A960").Value = $commonA
$ws.Range("B960").Value = $commonB
$ws.Range("C960").Value = $commonC
$ws.Range("D960").Value = 44585
$ws.Range("E960").Value = $commonE
$ws.Range("F960").Value = $commonF
$ws.Range("G960").Value = $commonG
$ws.Range("H960").Value = "Escarola"
$ws.Range("I960").Value = "Primera"
$ws.Range("J960").Value = 600
$ws.Range("K960").Value = 6000
$ws.Range("L960").Value = 6000
$ws.Range("M960").Value = 6000
$ws.Range("N960").Value = "$/caja 15 unidades"
$ws.Range("O960").Value = "Región del Maule"
$ws.Range("P960").Value = 400
$ws.Range("Q960").Value = 15
$ws.Range("R960").Value = $commonR

# Row 961 - Española, Primera
$ws.Range("A961").Value = $commonA
$ws.Range("B961").Value = $commonB
$ws.Range("C961").Value = $commonC
$ws.Range("D961").Value = 44585
$ws.Range("E961").Value = $commonE
$ws.Range("F961").Value = $commonF
$ws.Range("G961").Value = $commonG
$ws.Range("H961").Value = "Española"
$ws.Range("I961").Value = "Primera"
$ws.Range("J961").Value = 500
$ws.Range("K961").Value = 4500
$ws.Range("L961").Value = 4500
$ws.Range("M961").Value = 4500
$ws.Range("N961").Value = "$/caja 18 unidades"
$ws.Range("O961").Value = "Región del Maule"
$ws.Range("P961").Value = 250
$ws.Range("Q961").Value = 18
$ws.Range("R961").Value = $commonR

# Row 962 - Marina, Primera
$ws.Range("A962").Value = $commonA
$ws.Range("B962").Value = $commonB
$ws.Range("C962").Value = $commonC
$ws.Range("D962").Value = 44585
$ws.Range("E962").Value = $commonE
$ws.Range("F962").Value = $commonF
$ws.Range("G962").Value = $commonG
$ws.Range("H962").Value = "Marina"
$ws.Range("I962").Value = "Primera"
$ws.Range("J962").Value = 500
$ws.Range("K962").Value = 4500
$ws.Range("L962").Value = 4500
$ws.Range("M962").Value = 4500
$ws.Range("N962").Value = "$/caja 18 unidades"
$ws.Range("O962").Value = "Región del Maule"
$ws.Range("P962").Value = 250
$ws.Range("Q962").Value = 18
$ws.Range("R962").Value = $commonR
